$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows.
# D-column values that look purely numeric need NumberFormat forced to
# Text ("@") first, otherwise Excel auto-converts them to numbers on
# assignment (the source data keeps these as text strings).

$ws.Range("D2").Value = '66.804.20'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '3.075.79'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.35'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.79'
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '3.074.31'
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").Value = '  -1.89%  '
$ws.Range("E10").Value = '  -1.54%  '
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("E12").Value = '  -3.32%  '
$ws.Range("E13").Value = '  -2.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.82'
$ws.Range("E14").Value = '  -3.75%  '
$ws.Range("E15").Value = '  -1.79%  '
$ws.Range("D16").Value = '3.587.71'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").Value = '66.785.40'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.99'
$ws.Range("E19").Value = '  +4.11%  '
$ws.Range("D20").Value = '3.076.64'
$ws.Range("E20").Value = '  -1.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '492.45'
$ws.Range("E21").Value = '  +3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.71'
$ws.Range("E22").Value = '  -3.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.688'
$ws.Range("E23").Value = '  -3.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.68'
$ws.Range("E25").Value = '  -5.72%  '
$ws.Range("E26").Value = '  -3.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.15'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.82'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.27'
$ws.Range("E30").Value = '  -4.83%  '
$ws.Range("E31").Value = '  -2.50%  '
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("D34").Value = '0.0₃0915'
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.950'
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("E37").Value = '  -4.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '47.16'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.123'
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("E40").Value = '  -5.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.300'
$ws.Range("E41").Value = '  -3.49%  '
$ws.Range("E42").Value = '  -4.67%  '
$ws.Range("D43").Value = '2.754.72'
$ws.Range("E43").Value = '  -3.58%  '
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '135.48'
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.70'
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("E51").Value = '  -1.98%  '

# Rows 46/47 swap places: "Bittensor" moves up to rank 46 and
# "dogwifhat" moves down to rank 47, both with refreshed price/volume.
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '367.67'
$ws.Range("E46").Value = '  -5.21%  '

$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.50'
$ws.Range("E47").Value = '  -3.95%  '
